# nuevos experimentos no convexos
# Updates the follower-restriction coefficients, the modified point,
# and the resulting vectors (bf / BF / alpha) for the new (non-convex)
# experiment run.

$wb = $excel.ActiveWorkbook

# Helper: assign a value while forcing it to be stored as TEXT (shared
# string), matching the source data (these numeric-looking values are
# text in the workbook, not real numbers), then drop the temporary
# Number Format again so no stray cell style is left behind.
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Restricciones_del_follower (sheet index 3) ---
$ws3 = $wb.Worksheets.Item(3)

Set-TextValue $ws3.Range("A2") "10.077049180327869 - 2x_1 + 0.0983606557377048y_1 - 1.0819672131147542y_2"
Set-TextValue $ws3.Range("B2") "-7.57704918032787"
Set-TextValue $ws3.Range("D2") "0.96"
Set-TextValue $ws3.Range("E2") "0"
Set-TextValue $ws3.Range("F2") "0"

Set-TextValue $ws3.Range("A3") "-0.6872950819672126 + x_1 - 3x_2 - 0.0901639344262295y_1 + 0.9918032786885246y_2"
Set-TextValue $ws3.Range("B3") "-1.3127049180327874"
Set-TextValue $ws3.Range("D3") "0.7"
Set-TextValue $ws3.Range("E3") "0"
Set-TextValue $ws3.Range("F3") "0"

Set-TextValue $ws3.Range("A4") "-6.31 + x_1 + x_2"
Set-TextValue $ws3.Range("B4") "4.1"
Set-TextValue $ws3.Range("D4") "0.8"
Set-TextValue $ws3.Range("E4") "0"
Set-TextValue $ws3.Range("F4") "0"

# --- Punto_modificado (sheet index 4) ---
$ws4 = $wb.Worksheets.Item(4)

Set-TextValue $ws4.Range("A2") "4.5"
Set-TextValue $ws4.Range("B2") "1.6"
Set-TextValue $ws4.Range("C2") "6.1000000000000005"
Set-TextValue $ws4.Range("D2") "1.55"

# --- Vector_bf (sheet index 5) ---
$ws5 = $wb.Worksheets.Item(5)

Set-TextValue $ws5.Range("A2") "3.968688524590164"
Set-TextValue $ws5.Range("A3") "-0.6555737704918031"

# --- Vector_BF (sheet index 6; name collides case-insensitively with
#     "Vector_bf" so it is addressed by index, not by name) ---
$ws6 = $wb.Worksheets.Item(6)

Set-TextValue $ws6.Range("A2") "2.0"
Set-TextValue $ws6.Range("A3") "-1.0"
Set-TextValue $ws6.Range("A4") "-0.5"
Set-TextValue $ws6.Range("A5") "-0.0"

# --- Vector_Alpha (sheet index 7) -- these two cells are real numbers ---
$ws7 = $wb.Worksheets.Item(7)

$ws7.Range("A2").Value = 0.99
$ws7.Range("A3").Value = 0.09
